# Applies updated NOAA temperature data (average_county_temperature, column I)
# and the resulting recalculated worst/best ASHP COP values (columns N/O)
# for NAICS 311221 rows, per commit:
#   "Added merged + updated datasets / Updated temperature with NOAA data /
#    Added back NAICS 311230 / Added merged datasets"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 1.791666666666668
$ws.Range("N2").Value = 1.458486584262888
$ws.Range("O2").Value = 1.552746181345467
$ws.Range("I5").Value = 12.93898809523811
$ws.Range("N5").Value = 1.586442583591966
$ws.Range("O5").Value = 1.700608911205746
$ws.Range("I6").Value = 12.93898809523811
$ws.Range("I7").Value = 19.79629629629628
$ws.Range("N7").Value = 1.676945000770297
$ws.Range("O7").Value = 1.806427491177953
$ws.Range("I8").Value = 19.79629629629628
$ws.Range("I9").Value = 19.79629629629628
$ws.Range("N9").Value = 1.676945000770297
$ws.Range("O9").Value = 1.806427491177953
$ws.Range("I10").Value = 12.93898809523811
$ws.Range("I11").Value = 12.93898809523811
$ws.Range("N11").Value = 1.586442583591966
$ws.Range("O11").Value = 1.700608911205746
$ws.Range("I12").Value = 12.93898809523811
$ws.Range("I14").Value = 14.47727272727272
$ws.Range("N14").Value = 1.605884483070795
$ws.Range("O14").Value = 1.723253983867794
$ws.Range("I15").Value = 13.76976495726495
$ws.Range("N15").Value = 1.596883662077925
$ws.Range("O15").Value = 1.712764324418727
$ws.Range("I16").Value = 13.76976495726495
$ws.Range("I19").Value = 14.47727272727272
$ws.Range("N19").Value = 1.605884483070795
$ws.Range("O19").Value = 1.723253983867794
$ws.Range("I20").Value = 19.60879629629628
$ws.Range("N20").Value = 1.674333288469303
$ws.Range("O20").Value = 1.803359265239363
$ws.Range("I21").Value = 19.60879629629628
$ws.Range("I23").Value = 0.2777777777777778
$ws.Range("N23").Value = 1.442683896620278
$ws.Range("O23").Value = 1.534625267665953
$ws.Range("I25").Value = 14.47727272727272
$ws.Range("N25").Value = 1.605884483070795
$ws.Range("O25").Value = 1.723253983867794
$ws.Range("I26").Value = 14.47727272727272
$ws.Range("I30").Value = 12.93898809523811
$ws.Range("N30").Value = 1.586442583591966
$ws.Range("O30").Value = 1.700608911205746
$ws.Range("I31").Value = 12.93898809523811
$ws.Range("I32").Value = 13.76976495726495
$ws.Range("N32").Value = 1.596883662077925
$ws.Range("O32").Value = 1.712764324418727
$ws.Range("I33").Value = 14.47727272727272
$ws.Range("N33").Value = 1.605884483070795
$ws.Range("O33").Value = 1.723253983867794
$ws.Range("I34").Value = 14.47727272727272
$ws.Range("I35").Value = 5.462962962962945
$ws.Range("N35").Value = 1.49828630419821
$ws.Range("O35").Value = 1.598520446096654
$ws.Range("I36").Value = 5.462962962962945
$ws.Range("I37").Value = 19.60879629629628
$ws.Range("N37").Value = 1.674333288469303
$ws.Range("O37").Value = 1.803359265239363
$ws.Range("I38").Value = 19.60879629629628
$ws.Range("I39").Value = 14.47727272727272
$ws.Range("N39").Value = 1.605884483070795
$ws.Range("O39").Value = 1.723253983867794
$ws.Range("I40").Value = 21.28240740740739
$ws.Range("N40").Value = 1.697937058846468
$ws.Range("O40").Value = 1.831120384959332
$ws.Range("I41").Value = 14.47727272727272
$ws.Range("N41").Value = 1.605884483070795
$ws.Range("O41").Value = 1.723253983867794
